$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Data): 2023-06-23 -> 2023-06-26
$ws.Range("A2").Value = 45103

# Column B (Motivo): "confirmar resultado" -> "hemolise"
$ws.Range("B2").Value = "hemolise"

# Column C (Pedido): "789" -> "25899898-00 pd789878"
$ws.Range("C2").Value = "25899898-00 pd789878"

# Column D (Atendimento): "654" -> "588889"
# "588889" is purely numeric, so force a text format first so Excel keeps
# it as text instead of converting it to a number, then restore the
# default (unstyled) cell format so no stray formatting is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "588889"
$ws.Range("D2").Style = "Normal"
